$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.940.46'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.791.47'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '359.37'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '110.08'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.83%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.564'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.592'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.91'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.40'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.59'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.228.77'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.791.45'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.946'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +5.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.888.40'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.50'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.10'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.08'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '271.52'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.18'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.76'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.17%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.54'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.167'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +19.19%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.22'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('E30').Value = '  -1.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '52.36'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0467'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.88%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '33.96'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.76'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0845'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.24'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.22'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.33'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.55'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.88%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.115'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.24'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '120.24'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.16%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.90'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -10.15%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.086.62'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.24'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.14%  '
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.71'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.948'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.77%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.91'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.30%  '
